# Apply cell-value updates to Sheet1 as described by the target diff.
# The diff only changes numeric <v> contents of existing cells (no new
# rows/columns, no structural changes), so we just set .Value on each
# affected cell using the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("G62").Value = 0.104
$ws.Range("I62").Value = 0.005000000000000004
$ws.Range("G65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("G68").Value = 0.02200000000000002
$ws.Range("J68").Value = 0.01899999999999996
$ws.Range("E69").Value = -0.01100000000000001
$ws.Range("F69").Value = -0.02400000000000002
$ws.Range("E77").Value = 0.02000000000000002
$ws.Range("G77").Value = 0.04199999999999993
$ws.Range("I81").Value = 0.07699999999999996
$ws.Range("J81").Value = 0.13
$ws.Range("I84").Value = 0.03599999999999992
$ws.Range("J84").Value = 0.153
$ws.Range("G88").Value = 0.03900000000000003
$ws.Range("E89").Value = -0.03600000000000003
$ws.Range("J89").Value = -0.02000000000000002
$ws.Range("F90").Value = 0.001000000000000001
$ws.Range("K90").Value = -0.08500000000000002
$ws.Range("F91").Value = -0.03100000000000003
$ws.Range("H91").Value = 0.01999999999999996
$ws.Range("L93").Value = -0.04499999999999998
$ws.Range("K97").Value = -0.06999999999999995
$ws.Range("L97").Value = -0.05099999999999999
$ws.Range("E99").Value = -0.03400000000000003
$ws.Range("I99").Value = -0.03900000000000003
$ws.Range("F100").Value = 0.03899999999999998
$ws.Range("K100").Value = -0.106
